$d = $word.ActiveDocument

# The four inline pictures in this document's headers/footers each carry
# their filename twice in the underlying OOXML - once on <wp:docPr> and
# once on the nested <pic:cNvPr> - and both need to be kept in sync.
# Word's InlineShape object has no writable Name property (only the
# floating Shape object does, and even that only ever touches
# <wp:docPr>), so the rename is done by editing the package XML directly
# via Document.WordOpenXML, which round-trips the whole document (every
# part, with its own correctly-scoped relationship ids) losslessly.

$xml = $d.WordOpenXML

$singleline = [System.Text.RegularExpressions.RegexOptions]::Singleline

function Rename-PearsonLogo {
    param($xmlIn, [int]$id, [string]$descr, [string]$oldName, [string]$newName)

    $descrPat = [regex]::Escape($descr)
    $oldPat = [regex]::Escape($oldName)

    $pattern = '<wp:docPr descr="' + $descrPat + '" id="' + $id + '" name="' + $oldPat + '"/>(.*?)<pic:cNvPr descr="' + $descrPat + '" id="0" name="' + $oldPat + '"/>'
    $replacement = '<wp:docPr descr="' + $descr + '" id="' + $id + '" name="' + $newName + '"/>$1<pic:cNvPr descr="' + $descr + '" id="0" name="' + $newName + '"/>'

    return [regex]::Replace($xmlIn, $pattern, $replacement, $singleline)
}

$pearsonDescr = "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png"
$btecDescr = "BTec_Logo-Orange"

$xml = Rename-PearsonLogo $xml 2 $pearsonDescr "image1.png" "image2.png"
$xml = Rename-PearsonLogo $xml 4 $pearsonDescr "image1.png" "image2.png"
$xml = Rename-PearsonLogo $xml 1 $btecDescr "image2.jpg" "image1.jpg"
$xml = Rename-PearsonLogo $xml 3 $btecDescr "image2.jpg" "image1.jpg"

$d.WordOpenXML = $xml

Write-Host "done"
